$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 6: append ",R52" to the Designator for R1-10
$ws.Range("B6").Value = "R1-10,R52"

# Rows 9-10 currently carry blank placeholder styling; copy the formatting
# used by the other populated data rows (e.g. row 8) before filling values.
$ws.Range("A8:D8").Copy($ws.Range("A9:D9"))
$ws.Range("A8:D8").Copy($ws.Range("A10:D10"))

# Row 9: fill in previously-empty resistor row (5K, R51, R53)
$ws.Range("A9").Value = "5K"
$ws.Range("B9").Value = "R51, R53"
$ws.Range("C9").Value = 1206

# Row 10: new transistor row
$ws.Range("D10").Value = "C8545"
$ws.Range("A10").Value = "Transistor"
$ws.Range("B10").Value = "Q1"
$ws.Range("C10").Value = "SOT-23-3"

# Finish row 9
$ws.Range("D9").Value = "C17936"

# Update the selected cell shown in the saved view
$ws.Range("D16").Select()
